$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Duplicate the "Sprint 1" block (rows 5-9) as a new "Sprint 2" block
# (rows 12-15), separated by a blank thick-bottom-border row (row 11),
# exactly mirroring the existing layout/formatting.
# ---------------------------------------------------------------------------

# Copy all formatting (styles/number formats/borders) from the first sprint
# table down to where the second sprint table will live.
$ws.Range("A5:N9").Copy()
$ws.Range("A12").PasteSpecial(-4122)  # xlPasteFormats

# Recreate the merged-cell layout for the new block (mirrors A6:A8, etc.)
$ws.Range("E12:G12").Merge()
$ws.Range("H12:J12").Merge()

$ws.Range("A13:A15").Merge()
$ws.Range("B13:B15").Merge()
$ws.Range("C13:C15").Merge()
$ws.Range("D13:D15").Merge()
$ws.Range("E13:G15").Merge()
$ws.Range("H13:J13").Merge()
$ws.Range("H14:J14").Merge()
$ws.Range("H15:J15").Merge()
$ws.Range("K13:K15").Merge()
$ws.Range("L13:L15").Merge()
$ws.Range("M13:M15").Merge()
$ws.Range("N13:N15").Merge()

# Header row (row 12) repeats the same captions as row 5.
$ws.Range("A12").Value = "SPRINT"
$ws.Range("B12").Value = "INICIO"
$ws.Range("C12").Value = "DURACION"
$ws.Range("D12").Value = "Backlog ID"
$ws.Range("E12").Value = "Enunciado"
$ws.Range("H12").Value = "Tareas asociadas"
$ws.Range("K12").Value = "Tipo"
$ws.Range("L12").Value = "Estado"
$ws.Range("M12").Value = "Responsable"
$ws.Range("N12").Value = "Esfuerzo"

# Sprint 2 data.
$ws.Range("A13").Value = 2
$ws.Range("B13").Value = 45306
$ws.Range("C13").Value = "5 dias"
$ws.Range("D13").Value = "HT_02"
$ws.Range("E13").Value = "Como desarrollador, requiero utilizar PostgreSQL como sistema de gestión de base de datos para el sistema."
$ws.Range("H13").Value = "Configuración de PostgreSQL."
$ws.Range("H14").Value = "Diseño de esquemas de base de datos."
$ws.Range("H15").Value = "Implementación de operaciones de lectura y escritura."
$ws.Range("K13").Value = "Base de Datos"
$ws.Range("L13").Value = "Terminado"
$ws.Range("M13").Value = "Liliana Nogales"
$ws.Range("N13").Value = 12

# Row 15 (bottom of the table) carries the same thick-bottom border as row 8;
# give it the same "short" height Excel assigns a thick-border row lacking an
# explicit custom height (matches row 4's separator row).
$ws.Rows("15").RowHeight = 15.75

# A lone formatted (underlined) placeholder cell, left over after some
# scratch work further down the sheet.
$ws.Range("J22").Font.Underline = 2

# ---------------------------------------------------------------------------
# Misc view/page bits
# ---------------------------------------------------------------------------
$ws.Range("J22").Select()
$ws.PageSetup.PaperSize = 9      # xlPaperA4
$ws.PageSetup.Orientation = 1    # xlPortrait
